# Auto update: 2025-12-03 08:54:10
# Re-sync the DECISION sheet: the Bitcoin (BTC-USD) row and the Riot
# Platforms (RIOT) row swap places (row 2 <-> row 3 identity), and every
# ticker's scored metrics (종가/RSI/5일수익률/점수(룰)/확률 columns/최종점수)
# are refreshed with the latest run's numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes the Bitcoin USD / BTC-USD row ---
$ws.Range("B2").Value = "Bitcoin USD"
$ws.Range("C2").Value = "BTC-USD"
$ws.Range("D2").Value = 91451.64999999999
$ws.Range("E2").Value = 47.2
$ws.Range("F2").Value = 0.18
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = 46
$ws.Range("K2").Value = 50.8
$ws.Range("N2").Value = 66.04328690552585

# --- Row 3 becomes the Riot Platforms, Inc. / RIOT row ---
$ws.Range("B3").Value = "Riot Platforms, Inc."
$ws.Range("C3").Value = "RIOT"
$ws.Range("D3").Value = 15.22
$ws.Range("E3").Value = 44.3
$ws.Range("F3").Value = 9.65
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 43
$ws.Range("I3").Value = 43
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 49
$ws.Range("N3").Value = 66.04328690552585

# --- Row 4: Coinbase Global, Inc. / COIN (identity unchanged, metrics refreshed) ---
$ws.Range("D4").Value = 263.26
$ws.Range("E4").Value = 33.3
$ws.Range("F4").Value = 2.85
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 48.8
$ws.Range("N4").Value = 66.04328690552585

# --- Row 5: MARA Holdings, Inc. / MARA (identity unchanged, metrics refreshed) ---
$ws.Range("D5").Value = 11.91
$ws.Range("E5").Value = 32.8
$ws.Range("F5").Value = 6.24
$ws.Range("H5").Value = 43
$ws.Range("K5").Value = 48.8
$ws.Range("N5").Value = 66.04328690552585

# --- Row 6: Strategy Inc / MSTR (identity unchanged, metrics refreshed) ---
$ws.Range("D6").Value = 181.33
$ws.Range("E6").Value = 29.1
$ws.Range("F6").Value = 1.28
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 43
$ws.Range("K6").Value = 44.8
$ws.Range("N6").Value = 66.04328690552585
